$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price/volume snapshot (GitHub Actions nightly
# scrape). Every row's Price (D) / Volume(1h) (E) is updated to the new
# reading; rows 48-49 additionally swap rank position (Cronos <-> EnergySwap
# traded places), so Coin (B) and Link (C) are rewritten there too.
#
# The Price column holds plain text in the source data (coin values can
# look like "25.890.60" or "0.0\u2085...", not valid Excel numbers), so for
# every Price cell we briefly force Text number format before writing the
# string -- otherwise Excel would auto-coerce plain decimal-looking values
# (e.g. "215.76") into a real number and drop precision/trailing zeros --
# then restore the cell's default style so no formatting changes linger.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.890.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.640.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5080"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06460"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.666.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.267"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.867.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5656"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7703"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.914.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.398"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.984"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.252"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.760"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1227"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.855"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.245"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04982"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.318"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.260"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.578"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.388"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9093"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.583"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5531"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.128.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01575"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9991"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.507"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8003"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₈108"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4237"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.690"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05042"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.003"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
